$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.481.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.988.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8195"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +74.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "253.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3405"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +18.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +16.30%  "
$ws.Range("E10").Value = "  +8.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8430"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08112"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "101.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.988.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.481"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "31.458.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.727"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.249.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.0000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.943"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +11.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.682"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1572"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +64.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  +5.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.195"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.567"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.69%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.356"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.564"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.337"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05205"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.215"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7505"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.802"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.54%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9986"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02000"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.83%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.940"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.55%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.632"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.63%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "78.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.02%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4668"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.38%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.065"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.16%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "106.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.50%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8551"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.98%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.494"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.98%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4278"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.27%  "
